# Testing XL Trail on a new branch
#
# Adds a new row (row 4) to Sheet1 mirroring the existing row 3 pattern,
# introducing three new shared strings ("Property:Test2", "SQLTEXT",
# "Test2342342"), and switches the active sheet/selection from
# ProcessResponse back to Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ProcessResponse")

# --- New data row on Sheet1 (row 4) ---------------------------------------
$ws1.Range("F4").Value = "Data"
$ws1.Range("H4").Value = "Property:Test2"
$ws1.Range("I4").Value = "SQLTEXT"
$ws1.Range("K4").Value = "Test2342342"

# Carry over the same cell formatting used on row 3 for the matching columns.
$ws1.Range("B3").Copy()
$ws1.Range("B4").PasteSpecial(-4122)

$ws1.Range("K3").Copy()
$ws1.Range("K4").PasteSpecial(-4122)

$ws1.Range("Q3").Copy()
$ws1.Range("Q4").PasteSpecial(-4122)

# Touch the remaining "blank" cells on the row so they persist as empty
# cells (matching J3/M3-style/P3/S3 placeholders already present on row 3).
$ws1.Range("J4").Borders.LineStyle = -4142
$ws1.Range("M4").Borders.LineStyle = -4142
$ws1.Range("P4").Borders.LineStyle = -4142
$ws1.Range("S4").Borders.LineStyle = -4142

$ws1.Rows.Item(4).RowHeight = 18.75

# --- Selection / active sheet changes -------------------------------------
# ProcessResponse's selection moves to I3, and it is no longer the active tab.
$ws2.Activate()
$ws2.Range("I3").Select()

# Sheet1 becomes the active tab again, selection moves to M4.
$ws1.Activate()
$ws1.Range("M4").Select()
